# TC09_CDS_Filter_InstrumentModel-Illumina NextSeq.xlsx
# "Failed test cases from Instrument model filter"
#
# The Participant-tab Cypher query (row 2, column B -> "query") is rewritten
# to use OPTIONAL MATCH / apoc.coll.sort instead of the old plain MATCH, and
# the selected cell moves from D3 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lines = @(
    "MATCH (p:participant)-->(s:study)",
    "OPTIONAL MATCH (samp:sample)-->(p)",
    "OPTIONAL MATCH (p)<--(diag:diagnosis)",
    "OPTIONAL MATCH (samp)<--(f:file)",
    "OPTIONAL MATCH (f)<--(g:genomic_info)",
    "WITH s, p, samp, f, g, diag",
    "WHERE g.instrument_model in ['Illumina NextSeq']",
    "WITH p",
    "OPTIONAL MATCH (p)-->(s:study)",
    "OPTIONAL MATCH (samp:sample)-->(p)",
    "WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp",
    "RETURN ",
    "coalesce(p.participant_id,'') as ``Participant ID``,",
    "coalesce(s.study_name, '') as ``Study Name``,",
    "coalesce(s.phs_accession,'') as ``Accession``,",
    "coalesce(p.gender,'') as ``Gender``,",
    "coalesce(apoc.text.join(samp, ','), '') as ``Samples``",
    "ORDER BY p.participant_id limit 100"
)
$newParticipantQuery = [string]::Join("`r`n", $lines)

$ws.Range("B2").Value = $newParticipantQuery

# Move the active selection from D3 to C2.
$ws.Range("C2").Select()
